# Fixed naive component forecaster bug - Presentation state 11.02.
# Updates the forecast error table (ME, MAE, MSE, RMSE, SE, N) for Q0-Q9
# rows on the active worksheet to reflect the corrected naive forecaster
# output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = -0.03276308475587503
$ws.Range("C2").Value2 = 3.590912207024844
$ws.Range("D2").Value2 = 47.40170072604877
$ws.Range("E2").Value2 = 6.884889303834068
$ws.Range("F2").Value2 = 6.953318621525447
$ws.Range("G2").Value2 = 51

$ws.Range("B3").Value2 = 0.05364334786105588
$ws.Range("C3").Value2 = 3.574536985163549
$ws.Range("D3").Value2 = 44.07793895074698
$ws.Range("E3").Value2 = 6.639121850873576
$ws.Range("F3").Value2 = 6.706306911391149
$ws.Range("G3").Value2 = 50

$ws.Range("B4").Value2 = -0.04985138525762726
$ws.Range("C4").Value2 = 3.773358965575171
$ws.Range("D4").Value2 = 43.3640935270437
$ws.Range("E4").Value2 = 6.585141875999613
$ws.Range("F4").Value2 = 6.653192858124505
$ws.Range("G4").Value2 = 49

$ws.Range("B5").Value2 = 0.2063883000833151
$ws.Range("C5").Value2 = 3.726029611715911
$ws.Range("D5").Value2 = 41.36251091033633
$ws.Range("E5").Value2 = 6.431369287355246
$ws.Range("F5").Value2 = 6.496080503111255
$ws.Range("G5").Value2 = 48

$ws.Range("B6").Value2 = -0.03263097726318537
$ws.Range("C6").Value2 = 3.889169391173965
$ws.Range("D6").Value2 = 44.20905250212165
$ws.Range("E6").Value2 = 6.648988833057373
$ws.Range("F6").Value2 = 6.720790945475113
$ws.Range("G6").Value2 = 47

$ws.Range("B7").Value2 = 0.2010956551540337
$ws.Range("C7").Value2 = 3.91100830380399
$ws.Range("D7").Value2 = 49.94279819706635
$ws.Range("E7").Value2 = 7.06702187608517
$ws.Range("F7").Value2 = 7.142219545754433
$ws.Range("G7").Value2 = 46

$ws.Range("B8").Value2 = -0.1796802820817339
$ws.Range("C8").Value2 = 3.730114659890623
$ws.Range("D8").Value2 = 39.99999453259167
$ws.Range("E8").Value2 = 6.324554888100163
$ws.Range("F8").Value2 = 6.393439337599808
$ws.Range("G8").Value2 = 45

$ws.Range("B9").Value2 = 0.01305864851553937
$ws.Range("C9").Value2 = 3.680336663577857
$ws.Range("D9").Value2 = 42.37372025704332
$ws.Range("E9").Value2 = 6.509509985939289
$ws.Range("F9").Value2 = 6.584753686948571
$ws.Range("G9").Value2 = 44

$ws.Range("B10").Value2 = -0.1400871512439874
$ws.Range("C10").Value2 = 4.022154036027911
$ws.Range("D10").Value2 = 46.29094780507954
$ws.Range("E10").Value2 = 6.803745130814318
$ws.Range("F10").Value2 = 6.882806227928441
$ws.Range("G10").Value2 = 43

$ws.Range("B11").Value2 = 0.1042633263670164
$ws.Range("C11").Value2 = 3.990286954105756
$ws.Range("D11").Value2 = 47.50789127154154
$ws.Range("E11").Value2 = 6.892596845278384
$ws.Range("F11").Value2 = 6.975348331320163
$ws.Range("G11").Value2 = 42
